$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BADcluster_finalList")

# The "no-fill" style (style index 4 in cellXfs, used only by A13/A14) is
# redundant with the base "no-fill" style (index 0) - drop it so the style
# table collapses back down, same as a plain re-save would do.
$ws.Range("A13").ClearFormats()
$ws.Range("A14").ClearFormats()

# Remove the excluded subcluster "calcarine-astrocyte-3" row (row 2) -
# shifts all subsequent rows up by one.
$ws.Rows.Item(2).Delete()
